$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "harvester" column (B) for all data rows: Retrofitted_0684 -> S.GISH
$ws.Range("B2:B19").Value = "S.GISH"

# Add "experimentDesign" column (D) for all data rows: 90minuteInduction
$ws.Range("D2:D19").Value = "90minuteInduction"
$ws.Range("D2:D19").Font.Color = 0

# Select B19 to mirror the saved selection state in the target file
$ws.Range("B19").Select() | Out-Null
